$wb = $excel.ActiveWorkbook

# =========================================================================
# 1) Summary sheet: fix the "25.05.1017" typo -> "23.05.1017", add the new
#    "2017.05.30" row (row 8), and move the selection to C8.
# =========================================================================
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("A7").Value = "23.05.1017"

$summary.Range("A8").Value2 = 42885
$summary.Range("A8").NumberFormat = "m/d/yy"
$summary.Range("B8").Value = "Tutor list"

$summary.Range("C8").Select()

# =========================================================================
# 2) 2017.05.23 sheet: tutors worked hours got filled in -> update column F
#    (Worked) for the first, second and fourth data rows; formulas in G and
#    the totals row recalc automatically. Move the selection to I11.
# =========================================================================
$week0523 = $wb.Worksheets.Item("2017.05.23")

$week0523.Range("F3").Value = 4
$week0523.Range("F4").Value = 0
$week0523.Range("F7").Value = 6

$week0523.Range("I11").Select()

# =========================================================================
# 3) Brand new weekly sheet "2017.05.30", appended after "2017.05.23".
#    Layout mirrors the other weekly sheets, but now also stores each
#    tutor's subject (column B) next to their name (column A).
# =========================================================================
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$week0530 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$week0530.Name = "2017.05.30"

$week0530.Columns.Item(2).ColumnWidth = 15.45

# --- header row ---
$week0530.Range("A1").Value = "User Story"
$week0530.Range("B1").Value = "Task"
$week0530.Range("C1").Value = "Initial Estimation"
$week0530.Range("D1").Value = "Current Estimation"
$week0530.Range("E1").Value = "Worked"
$week0530.Range("F1").Value = "Remain"
$week0530.Range("J1").Value = "Name"
$week0530.Range("K1").Value = "Planned hours"

# --- row 2: Liste sortieren / Programmieren ---
$week0530.Range("A2").Value = "Liste sortieren"
$week0530.Range("B2").Value = "Programmieren"
$week0530.Range("C2").Value = 1
$week0530.Range("D2").Value = 1
$week0530.Range("F2").Formula = "=D2-E2"

# --- row 3: Liste filtern / Programmieren ---
$week0530.Range("A3").Value = "Liste filtern"
$week0530.Range("B3").Value = "Programmieren"
$week0530.Range("C3").Value = 3
$week0530.Range("D3").Value = 3
$week0530.Range("F3").Formula = "=D3-E3"

# --- rows 4-5: remain-only shared formula block ---
$week0530.Range("F4:F5").Formula = "=D4-E4"

# --- tutors / planned hours ---
$week0530.Range("J2").Value = "Eva"
$week0530.Range("K2").Value = 4
$week0530.Range("J3").Value = "Danijal"
$week0530.Range("K3").Value = 3
$week0530.Range("K4").Formula = "=SUM(K2:K3)"

# --- row 6: Tutor list / Bilder in Liste anzeigen ---
$week0530.Range("A6").Value = "Tutor list"
$week0530.Range("B6").Value = "Bilder in Liste anzeigen"
$week0530.Range("C6").Value = 3
$week0530.Range("D6").Value = 3
$week0530.Range("E6").Value = 0
$week0530.Range("F6").Formula = "=D6-E6"

# --- rows 7-19: remain-only shared formula block ---
$week0530.Range("F7:F19").Formula = "=D7-E7"

# --- totals row ---
$week0530.Range("C23").Formula = "=SUM(C2:C22)"
$week0530.Range("D23").Formula = "=SUM(D2:D22)"
$week0530.Range("E23").Formula = "=SUM(E2:E22)"
$week0530.Range("F23").Formula = "=SUM(F2:F22)"

$week0530.Range("D3").Select()
$week0530.Activate()
